$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new day column (BP) with the 20-aug prices.
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Header cell: copy the formatting of the previous header (BO1) so the new
# header keeps the bold / centered / bordered look, then set its text.
$wsSpot.Range("BO1").Copy()
$wsSpot.Range("BP1").PasteSpecial(-4122)
$wsSpot.Range("BP1").Value = "20-aug"

$spotValues = @(82.90000000000001, 74.12, 67.05, 66.08, 64.95, 70.59999999999999, `
    74.06999999999999, 78.43000000000001, 89.81, 71.55, 49.8, 56.03, 49.02, 35, `
    31.53, 27.2, 27.99, 37.72, 54.11, 69.27, 78.73, 100.07, 97.52, 84.98)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 68).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append the 2025-08-18 quote as a new row (65).
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A65").Value = "'2025-08-18"
$wsGaz.Range("A65").Style = "Normal"
$wsGaz.Range("B65").Value = 29.95

# ---------------------------------------------------------------------------
# Sheet "CO2": append the 2025-08-18 quote as a new row (65).
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A65").Value = "'2025-08-18"
$wsCo2.Range("A65").Style = "Normal"
$wsCo2.Range("B65").Value = 71
